$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-7) got reshuffled/re-valued in columns D and I:Q.
# Columns A,B,C,E,F,G,H,O,R are identical for every row and remain unchanged.
# New values per row (matching the target diff):
$rows = @{
    2 = @{ D = 45204; I = "Primera"; J = 100; K = 800;  L = 1000; M = 900;  N = '$/docena de matas'; P = 150;  Q = 6 }
    3 = @{ D = 44267; I = "Primera"; J = 120; K = 1500; L = 1800; M = 1650; N = '$/docena de matas'; P = 275;  Q = 6 }
    4 = @{ D = 44623; I = "Primera"; J = 300; K = 1800; L = 2000; M = 1900; N = '$/paquete';          P = 1900; Q = 1 }
    5 = @{ D = 44377; I = "Segunda"; J = 550; K = 2000; L = 2800; M = 2364; N = '$/docena de matas'; P = 394;  Q = 6 }
    6 = @{ D = 44370; I = "Segunda"; J = 100; K = 1000; L = 1200; M = 1080; N = '$/docena de matas'; P = 180;  Q = 6 }
    7 = @{ D = 45218; I = "Primera"; J = 180; K = 1400; L = 1500; M = 1444; N = '$/docena de matas'; P = 241;  Q = 6 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $ws.Range("D$r").Value = $vals.D
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
}
